# Updated TPM-derived NATMI ligand-receptor metrics for Plat -> Lrp1
# (re-run of the scoring pipeline with refreshed per-cluster TPM input).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.074329666666667
$ws.Range("H2").Value = 27.222989
$ws.Range("I2").Value = 0.1061942826727393
$ws.Range("J2").Value = 0.1061942826727393
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 17.41278866812122
$ws.Range("R2").Value = 156.715098013091
$ws.Range("S2").Value = 0.0006928135351058298
$ws.Range("T2").Value = 0.0006928135351058298

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.074329666666667
$ws.Range("H3").Value = 27.222989
$ws.Range("I3").Value = 0.1061942826727393
$ws.Range("J3").Value = 0.1061942826727393
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 1645.070627698343
$ws.Range("R3").Value = 14805.63564928509
$ws.Range("S3").Value = 0.06545345600851583
$ws.Range("T3").Value = 0.06545345600851583

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.074329666666667
$ws.Range("H4").Value = 27.222989
$ws.Range("I4").Value = 0.1061942826727393
$ws.Range("J4").Value = 0.1061942826727393
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 263.5881548019793
$ws.Range("R4").Value = 2372.293393217813
$ws.Range("S4").Value = 0.01048754710236116
$ws.Range("T4").Value = 0.01048754710236116

# Row 5: ECs -> Resolving-Mac
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.074329666666667
$ws.Range("H5").Value = 27.222989
$ws.Range("I5").Value = 0.1061942826727393
$ws.Range("J5").Value = 0.1061942826727393
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 742.9562526899249
$ws.Range("R5").Value = 6686.606274209325
$ws.Range("S5").Value = 0.02956046602675644
$ws.Range("T5").Value = 0.02956046602675644

# Row 6: FAPs -> ECs
$ws.Range("I6").Value = 0.6709510795864623
$ws.Range("J6").Value = 0.6709510795864622
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 110.0165570258707
$ws.Range("R6").Value = 990.1490132328361
$ws.Range("S6").Value = 0.004377297700328063
$ws.Range("T6").Value = 0.004377297700328063

# Row 7: FAPs -> FAPs
$ws.Range("I7").Value = 0.6709510795864623
$ws.Range("J7").Value = 0.6709510795864622
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.4135445512346047
$ws.Range("T7").Value = 0.4135445512346046

# Row 8: FAPs -> MuSCs
$ws.Range("I8").Value = 0.6709510795864623
$ws.Range("J8").Value = 0.6709510795864622
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 1665.388687408039
$ws.Range("R8").Value = 14988.49818667235
$ws.Range("S8").Value = 0.06626186338324822
$ws.Range("T8").Value = 0.06626186338324822

# Row 9: FAPs -> Resolving-Mac
$ws.Range("I9").Value = 0.6709510795864623
$ws.Range("J9").Value = 0.6709510795864622
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 4694.106756801723
$ws.Range("R9").Value = 42246.96081121551
$ws.Range("S9").Value = 0.1867673672682813
$ws.Range("T9").Value = 0.1867673672682813

# Row 10: MuSCs -> ECs
$ws.Range("G10").Value = 18.86476133333333
$ws.Range("H10").Value = 56.594284
$ws.Range("I10").Value = 0.2207689020760095
$ws.Range("J10").Value = 0.2207689020760095
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 36.19970999935511
$ws.Range("R10").Value = 325.797389994196
$ws.Range("S10").Value = 0.00144030054762992
$ws.Range("T10").Value = 0.00144030054762992

# Row 11: MuSCs -> FAPs
$ws.Range("G11").Value = 18.86476133333333
$ws.Range("H11").Value = 56.594284
$ws.Range("I11").Value = 0.2207689020760095
$ws.Range("J11").Value = 0.2207689020760095
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 3419.962235007269
$ws.Range("R11").Value = 30779.66011506542
$ws.Range("S11").Value = 0.1360721806899107
$ws.Range("T11").Value = 0.1360721806899107

# Row 12: MuSCs -> MuSCs
$ws.Range("G12").Value = 18.86476133333333
$ws.Range("H12").Value = 56.594284
$ws.Range("I12").Value = 0.2207689020760095
$ws.Range("J12").Value = 0.2207689020760095
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 547.9774058572032
$ws.Range("R12").Value = 4931.796652714828
$ws.Range("S12").Value = 0.02180272045712557
$ws.Range("T12").Value = 0.02180272045712557

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("G13").Value = 18.86476133333333
$ws.Range("H13").Value = 56.594284
$ws.Range("I13").Value = 0.2207689020760095
$ws.Range("J13").Value = 0.2207689020760095
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 1544.54300239806
$ws.Range("R13").Value = 13900.88702158254
$ws.Range("S13").Value = 0.06145370038134335
$ws.Range("T13").Value = 0.06145370038134335

# Row 14: Resolving-Mac -> ECs
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.1782266666666667
$ws.Range("H14").Value = 0.53468
$ws.Range("I14").Value = 0.002085735664789058
$ws.Range("J14").Value = 0.002085735664789059
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 0.3420002794355556
$ws.Range("R14").Value = 3.07800251492
$ws.Range("S14").Value = 0.0000136073794450119
$ws.Range("T14").Value = 0.0000136073794450119

# Row 15: Resolving-Mac -> FAPs
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.1782266666666667
$ws.Range("H15").Value = 0.53468
$ws.Range("I15").Value = 0.002085735664789058
$ws.Range("J15").Value = 0.002085735664789059
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 32.31042569270222
$ws.Range("R15").Value = 290.79383123432
$ws.Range("S15").Value = 0.001285555155557431
$ws.Range("T15").Value = 0.001285555155557432

# Row 16: Resolving-Mac -> MuSCs
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.1782266666666667
$ws.Range("H16").Value = 0.53468
$ws.Range("I16").Value = 0.002085735664789058
$ws.Range("J16").Value = 0.002085735664789059
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 5.177069814395556
$ws.Range("R16").Value = 46.59362832956
$ws.Range("S16").Value = 0.000205983321107409
$ws.Range("T16").Value = 0.0002059833211074091

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.1782266666666667
$ws.Range("H17").Value = 0.53468
$ws.Range("I17").Value = 0.002085735664789058
$ws.Range("J17").Value = 0.002085735664789059
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 14.59222017054222
$ws.Range("R17").Value = 131.32998153488
$ws.Range("S17").Value = 0.0005805898086792062
$ws.Range("T17").Value = 0.0005805898086792063
